# Append a "Database" Heading2 section plus six schema tables
# (users, Folders, Files, user_folder, user_file, folder_file) right
# after the final "Password" bullet, before the trailing empty paragraph.
$d = $word.ActiveDocument

# The document body's last paragraph is the trailing empty <w:p/> that sits
# immediately before <w:sectPr>. Replacing its range with our OOXML fragment
# inserts the new heading/tables there and re-creates that trailing empty
# paragraph as the fragment's final element.
$range = $d.Paragraphs.Last.Range

$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:t>Database</w:t></w:r></w:p><w:tbl><w:tblPr><w:tblStyle w:val="TableGrid"/><w:tblW w:w="0" w:type="auto"/><w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/></w:tblPr><w:tblGrid><w:gridCol w:w="1248"/><w:gridCol w:w="1659"/><w:gridCol w:w="1637"/><w:gridCol w:w="3039"/><w:gridCol w:w="1433"/></w:tblGrid><w:tr><w:tc><w:tcPr><w:tcW w:w="9016" w:type="dxa"/><w:gridSpan w:val="5"/></w:tcPr><w:p><w:r><w:t>users</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="1852" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>id</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1915" w:type="dxa"/></w:tcPr><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>first_name</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1904" w:type="dxa"/></w:tcPr><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>last_name</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1758" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>email</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1587" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>password</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="1852" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>1</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1915" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>Oluwatobi</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1904" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>Sofela</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1758" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>contact@codesweetly.com</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1587" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>test</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="1852" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>2</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1915" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>Sarah</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1904" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>Precious</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1758" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>sp@example.com</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1587" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>example</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="1852" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>3</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1915" w:type="dxa"/></w:tcPr><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>Dav</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1904" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>Emma</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1758" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>de@sample.com</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1587" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>sample</w:t></w:r></w:p></w:tc></w:tr></w:tbl><w:p/><w:tbl><w:tblPr><w:tblStyle w:val="TableGrid"/><w:tblW w:w="0" w:type="auto"/><w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/></w:tblPr><w:tblGrid><w:gridCol w:w="3005"/><w:gridCol w:w="3005"/><w:gridCol w:w="3006"/></w:tblGrid><w:tr><w:tc><w:tcPr><w:tcW w:w="9016" w:type="dxa"/><w:gridSpan w:val="3"/></w:tcPr><w:p><w:r><w:t>Folders</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="3005" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>id</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3005" w:type="dxa"/></w:tcPr><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>folder_name</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3006" w:type="dxa"/></w:tcPr><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>updated_time</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="3005" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>1</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3005" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>Flowers</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3006" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>2024-10-22</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="3005" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>2</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3005" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>Office Jokes</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3006" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>2024-10-24</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="3005" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>3</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3005" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>Articles</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3006" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>2024-10-28</w:t></w:r></w:p></w:tc></w:tr></w:tbl><w:p/><w:tbl><w:tblPr><w:tblStyle w:val="TableGrid"/><w:tblW w:w="0" w:type="auto"/><w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/></w:tblPr><w:tblGrid><w:gridCol w:w="2100"/><w:gridCol w:w="2340"/><w:gridCol w:w="2542"/><w:gridCol w:w="2034"/></w:tblGrid><w:tr><w:tc><w:tcPr><w:tcW w:w="9016" w:type="dxa"/><w:gridSpan w:val="4"/></w:tcPr><w:p><w:r><w:t>Files</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="2100" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>id</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2340" w:type="dxa"/></w:tcPr><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>file_data</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2542" w:type="dxa"/></w:tcPr><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>updated_time</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2034" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>in-folder</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="2100" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>1</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2340" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>object</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2542" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>2024-10-22</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2034" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>true</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="2100" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>2</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2340" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>object</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2542" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>2024-10-22</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2034" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>false</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="2100" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>3</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2340" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>Object</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2542" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>2024-10-25</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="2034" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>true</w:t></w:r></w:p></w:tc></w:tr></w:tbl><w:p/><w:tbl><w:tblPr><w:tblStyle w:val="TableGrid"/><w:tblW w:w="0" w:type="auto"/><w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/></w:tblPr><w:tblGrid><w:gridCol w:w="3005"/><w:gridCol w:w="3005"/><w:gridCol w:w="3006"/></w:tblGrid><w:tr><w:tc><w:tcPr><w:tcW w:w="9016" w:type="dxa"/><w:gridSpan w:val="3"/></w:tcPr><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>user_folder</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="3005" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>id</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3005" w:type="dxa"/></w:tcPr><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>user_id</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3006" w:type="dxa"/></w:tcPr><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>folder_id</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="3005" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>1</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3005" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>3</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3006" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>1</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="3005" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>2</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3005" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>1</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3006" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>3</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="3005" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>3</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3005" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>2</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3006" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>2</w:t></w:r></w:p></w:tc></w:tr></w:tbl><w:p/><w:tbl><w:tblPr><w:tblStyle w:val="TableGrid"/><w:tblW w:w="0" w:type="auto"/><w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/></w:tblPr><w:tblGrid><w:gridCol w:w="3005"/><w:gridCol w:w="3005"/><w:gridCol w:w="3006"/></w:tblGrid><w:tr><w:tc><w:tcPr><w:tcW w:w="9016" w:type="dxa"/><w:gridSpan w:val="3"/></w:tcPr><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>user_f</w:t></w:r><w:r><w:t>ile</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="3005" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>id</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3005" w:type="dxa"/></w:tcPr><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>user_id</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3006" w:type="dxa"/></w:tcPr><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>f</w:t></w:r><w:r><w:t>ile</w:t></w:r><w:r><w:t>_id</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="3005" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>1</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3005" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>2</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3006" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>2</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="3005" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>2</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3005" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>3</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3006" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>3</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="3005" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>3</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3005" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>2</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3006" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>1</w:t></w:r></w:p></w:tc></w:tr></w:tbl><w:p/><w:tbl><w:tblPr><w:tblStyle w:val="TableGrid"/><w:tblW w:w="0" w:type="auto"/><w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/></w:tblPr><w:tblGrid><w:gridCol w:w="3005"/><w:gridCol w:w="3005"/><w:gridCol w:w="3006"/></w:tblGrid><w:tr><w:tc><w:tcPr><w:tcW w:w="9016" w:type="dxa"/><w:gridSpan w:val="3"/></w:tcPr><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>folder_file</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="3005" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>id</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3005" w:type="dxa"/></w:tcPr><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>folder_id</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3006" w:type="dxa"/></w:tcPr><w:p><w:proofErr w:type="spellStart"/><w:r><w:t>file_id</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="3005" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>1</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3005" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>2</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3006" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>1</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="3005" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>2</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3005" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>1</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3006" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>3</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="3005" w:type="dxa"/></w:tcPr><w:p><w:r><w:lastRenderedPageBreak/><w:t>3</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3005" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>1</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="3006" w:type="dxa"/></w:tcPr><w:p><w:r><w:t>2</w:t></w:r></w:p></w:tc></w:tr></w:tbl><w:p/>'

$null = $range.InsertXML($xml)
